$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Widen columns D and F on the "CUMPLIMIENTO MENSUAL" sheet
$ws.Columns.Item(4).ColumnWidth = 12.17
$ws.Columns.Item(6).ColumnWidth = 25.17

# Row 17 - PUERTAS DE SEGURIDAD
$ws.Range("D17").Value = -124.22
$ws.Range("E17").Value = 466.22
$ws.Range("F17").Value = -0.3632163742690058

# Row 18 - SAL SOLUBLE
$ws.Range("D18").Value = -248.83
$ws.Range("E18").Value = 3048.83
$ws.Range("F18").Value = -0.08886785714285715

# Row 19 - TOTAL
$ws.Range("D19").Value = -373.05
$ws.Range("E19").Value = 50760.24762291768
$ws.Range("F19").Value = -0.007403666359693026
